# Update "想去人数" (interest count) figures in the "展览" and "全部类型"
# sheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 62
$ws1.Range("F5").Value = 574
$ws1.Range("F12").Value = 1062
$ws1.Range("F14").Value = 12197
$ws1.Range("F15").Value = 12708
$ws1.Range("F22").Value = 20

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 62
$ws4.Range("F6").Value = 574
$ws4.Range("F13").Value = 1062
$ws4.Range("F15").Value = 12197
$ws4.Range("F16").Value = 12708
$ws4.Range("F23").Value = 20
